$d = $word.ActiveDocument

# The first two paragraphs:
#   "Reflections During Advent, Part Three""Chastity ====..."
#   "By Dorothy Day"
# become a single Title-styled paragraph "Dorothy Day" (pandoc-style
# title block), split into three runs: "Dorothy", " ", "Day".

# Drop the first paragraph (with its trailing paragraph mark) entirely -
# this merges what was paragraph two up to become paragraph one.
$d.Paragraphs(1).Range.Delete()

$titlePara = $d.Paragraphs(1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p>' +
       '<w:pPr><w:pStyle w:val="Title"/></w:pPr>' +
       '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' +
       '</w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

# Replace the (now plain "By Dorothy Day") paragraph's contents with the
# Title-styled, three-run "Dorothy Day" paragraph.
$titlePara.Range.InsertXML($xml)
